$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number-format on D2:E51 before writing, so numeric-looking
# strings (e.g. "539.95") are stored as text (matching the original
# inlineStr cells) instead of being auto-coerced into floats.
$fmtRange = $ws.Range("D2:E51")
$fmtRange.NumberFormat = "@"

$ws.Range("D2").Value = "58.941.96"
$ws.Range("E2").Value = "  -3.15%  "
$ws.Range("D3").Value = "3.233.65"
$ws.Range("E3").Value = "  -3.68%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "539.95"
$ws.Range("E5").Value = "  -4.68%  "
$ws.Range("D6").Value = "136.66"
$ws.Range("E6").Value = "  -7.86%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "3.232.88"
$ws.Range("E8").Value = "  -3.64%  "
$ws.Range("D9").Value = "0.460"
$ws.Range("E9").Value = "  -4.18%  "
$ws.Range("D10").Value = "7.64"
$ws.Range("E10").Value = "  -3.90%  "
$ws.Range("E11").Value = "  -5.49%  "
$ws.Range("D12").Value = "0.395"
$ws.Range("E12").Value = "  -4.24%  "
$ws.Range("D13").Value = "3.787.03"
$ws.Range("E13").Value = "  -3.71%  "
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("D15").Value = "26.01"
$ws.Range("E15").Value = "  -6.96%  "
$ws.Range("D16").Value = "3.236.21"
$ws.Range("E16").Value = "  -3.69%  "
$ws.Range("E17").Value = "  -5.64%  "
$ws.Range("D18").Value = "58.962.96"
$ws.Range("E18").Value = "  -3.31%  "
$ws.Range("D19").Value = "5.90"
$ws.Range("E19").Value = "  -6.98%  "
$ws.Range("E20").Value = "  -6.00%  "
$ws.Range("E21").Value = "  -6.19%  "
$ws.Range("D22").Value = "362.05"
$ws.Range("E22").Value = "  -3.12%  "
$ws.Range("D24").Value = "70.62"
$ws.Range("E24").Value = "  -6.35%  "
$ws.Range("E25").Value = "  -6.84%  "
$ws.Range("D26").Value = "3.367.40"
$ws.Range("E26").Value = "  -3.87%  "
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").Value = "0.0₃0977"
$ws.Range("E27").Value = "  -9.79%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.170"
$ws.Range("E28").Value = "  -3.06%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "7.09"
$ws.Range("E30").Value = "  -3.54%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  -6.41%  "
$ws.Range("D33").Value = "7.12"
$ws.Range("E33").Value = "  -7.12%  "
$ws.Range("E34").Value = "  -3.74%  "
$ws.Range("E35").Value = "  -4.72%  "
$ws.Range("D36").Value = "4.96"
$ws.Range("E36").Value = "  -7.28%  "
$ws.Range("D37").Value = "163.39"
$ws.Range("E37").Value = "  -3.07%  "
$ws.Range("D38").Value = "6.43"
$ws.Range("E38").Value = "  -4.78%  "
$ws.Range("E39").Value = "  -6.55%  "
$ws.Range("D40").Value = "26.43"
$ws.Range("E40").Value = "  -9.62%  "
$ws.Range("D41").Value = "0.0713"
$ws.Range("E41").Value = "  -4.88%  "
$ws.Range("D42").Value = "3.265.81"
$ws.Range("E42").Value = "  -3.83%  "
$ws.Range("D43").Value = "41.21"
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("D44").Value = "0.718"
$ws.Range("E44").Value = "  -5.45%  "
$ws.Range("D45").Value = "1.11"
$ws.Range("E45").Value = "  -2.84%  "
$ws.Range("E46").Value = "  -5.73%  "
$ws.Range("E47").Value = "  -5.98%  "
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("D49").Value = "2.300.70"
$ws.Range("E49").Value = "  -7.89%  "
$ws.Range("D50").Value = "6.31"
$ws.Range("E50").Value = "  -5.43%  "
$ws.Range("D51").Value = "20.97"
$ws.Range("E51").Value = "  -7.03%  "

# Restore default (General) styling so we don't leave a stray
# text-format style attached to these cells.
$fmtRange.Style = "Normal"

